# ---------------------------------------------------------------------------
# Commit: "Wed, Jul 08, 2020  4:05:01 AM"
#
# 1) The table on slide 16 (the cash-flow glossary table) gets re-styled:
#    its custom "Table_0" style ({57EBF60A-993D-4A4F-BEE1-366F9E5FD68E}) is
#    swapped for the built-in style {53F2B751-ED73-47C1-A6DF-AE8AE5F9BFD7}.
#
# 2) The deck's theme colour palette is swapped from the custom "Integral"
#    palette back to the stock Office palette (the notes master keeps its
#    own separate theme part, which the PowerPoint object model does not
#    expose a handle to - only the slide master / presentation theme's
#    ThemeColorScheme is reachable here - so we repaint every one of its
#    12 colour slots to the stock "Office" RGB values).
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- 1) Table style on slide 16 -------------------------------------------
$slide = $p.Slides.Item(16)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shape = $slide.Shapes.Item($i)
    if ($shape.HasTable) {
        $shape.Table.ApplyStyle("{53F2B751-ED73-47C1-A6DF-AE8AE5F9BFD7}")
    }
}

# --- 2) Theme colours: Integral -> Office -----------------------------------
$theme = $p.SlideMaster.Theme
$colors = $theme.ThemeColorScheme

$colors.Colors(1).RGB  = 0          # dk1      000000
$colors.Colors(2).RGB  = 16777215   # lt1      FFFFFF
$colors.Colors(3).RGB  = 6968388    # dk2      44546A
$colors.Colors(4).RGB  = 15132391   # lt2      E7E6E6
$colors.Colors(5).RGB  = 13998939   # accent1  5B9BD5
$colors.Colors(6).RGB  = 3243501    # accent2  ED7D31
$colors.Colors(7).RGB  = 10855845   # accent3  A5A5A5
$colors.Colors(8).RGB  = 49407      # accent4  FFC000
$colors.Colors(9).RGB  = 12874308   # accent5  4472C4
$colors.Colors(10).RGB = 4697456    # accent6  70AD47
$colors.Colors(11).RGB = 12673797   # hlink    0563C1
$colors.Colors(12).RGB = 7491477    # folHlink 954F72
